$wb = $excel.ActiveWorkbook

# Rename the existing (only) sheet to "Regression"
$ws1 = $wb.ActiveSheet
$ws1.Name = "Regression"

# Add the new worksheets (added after the last sheet, in order)
$wsSmoke  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSmoke.Name = "Smoke"

$wsSanity = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSanity.Name = "Sanity"

$wsSheet4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSheet4.Name = "Sheet4"

$wsSheet5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSheet5.Name = "Sheet5"

$wsSheet6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSheet6.Name = "Sheet6"

# Update the selection on the Regression sheet
[void]$ws1.Range("H2").Select()

# Set the selection on the Sanity sheet
[void]$wsSanity.Range("E18").Select()

# Re-select the Regression sheet so it stays the active tab
[void]$ws1.Select()
